# "fixed ship highlighting" - the last simulated time step (row 9, time
# step 7) was still carrying "everything matches" placeholder data, and
# three bogus extra rows (10-12, time steps 8-10) were left over from the
# test harness. Correct row 9 so only the genuinely mismatched sensed
# values are sensed as off, then drop the stray trailing rows from both
# result sheets and shrink the conditional-formatting highlight range and
# sheet dimension to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "USS Test Ship 1"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 9 (time step 7): ship truth/sensed state + sys truth/sensed state
# all flip from 1 -> 0.
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 0
$ws1.Range("D9").Value = 0
$ws1.Range("E9").Value = 0

# Drop the extra rows 10-12 (time steps 8-10) entirely.
$ws1.Rows("10:12").Delete()

# The sheet no longer needs to be shown with a frozen header row.
$ws1.Activate()
$excel.ActiveWindow.FreezePanes = $false

# Shrink the color-scale highlight on column F to the new used range.
$ws1.Range("F2:F9").FormatConditions.Item(1).ModifyAppliesToRange($ws1.Range("F2:F9"))

# ---------------------------------------------------------------
# Sheet 2: "System 1 History"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Row 9 (time step 7): component/system truth & sensed values that drive
# the mismatch flip from 1 -> 0 (B, D, F, H); C/E/G/I are unaffected.
$ws2.Range("B9").Value = 0
$ws2.Range("D9").Value = 0
$ws2.Range("F9").Value = 0
$ws2.Range("H9").Value = 0

# Drop the extra rows 10-12 (time steps 8-10) entirely.
$ws2.Rows("10:12").Delete()

# Shrink the color-scale highlight on column J to the new used range.
$ws2.Range("J2:J9").FormatConditions.Item(1).ModifyAppliesToRange($ws2.Range("J2:J9"))
